$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 239, shifting existing rows 239:302 down to 240:303
$ws.Rows("239:239").Insert()

# Populate the newly inserted row 239 with the new data record
$ws.Range("A239").Value = 3
$ws.Range("B239").Value = "Femacal de La Calera"
$ws.Range("C239").Value = "Coquimbo"
$ws.Range("D239").NumberFormat = $ws.Range("D240").NumberFormat
$ws.Range("D239").Value = 44508
$ws.Range("E239").Value = 5
$ws.Range("F239").Value = 100112021
$ws.Range("G239").Value = "Ají"
$ws.Range("H239").Value = "Americana (o)"
$ws.Range("I239").Value = "Primera"
$ws.Range("J239").Value = 73
$ws.Range("K239").Value = 33000
$ws.Range("L239").Value = 35000
$ws.Range("M239").Value = 33959
$ws.Range("N239").Value = "`$/caja 15 kilos"
$ws.Range("O239").Value = "Limache"
$ws.Range("P239").Value = 2264
$ws.Range("Q239").Value = 15
$ws.Range("R239").Value = "Hortaliza"
